# Applies the "add edited checklist mobile testing file" commit:
#  - Numbers the ID column (A) on sheet "Checklist mobile app testing"
#    sequentially (1..122) for every checklist-item row.
#  - Adds print / page-setup configuration to that sheet (fit to page,
#    centred printing, gridlines, custom margins, A4 portrait).
#  - Re-themes the "Legal" heading font (was the non-installed
#    "SF Pro Display") to the workbook's Arial minor theme font.
#  - Normalises the font used by the "App Review" section's ID column
#    cells (A120:A128) to match the rest of the sheet.
#  - Removes the three trailing blank rows (37:39) that were left over
#    at the bottom of the "iOS App Store Review Guidelines" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Checklist mobile app testing")
$ws2 = $wb.Worksheets.Item("iOS App Store Review Guidelines")

# ---------------------------------------------------------------------
# 1. Fill in sequential checklist IDs in column A of sheet 1.
# ---------------------------------------------------------------------
$idRows = @(4,5,6,7,8,9,10,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,48,49,50,51,52,53,54,55,56,57,59,60,61,62,63,65,66,67,69,70,71,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,89,90,92,93,94,95,96,97,98,99,100,101,102,103,104,105,107,108,109,110,111,112,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,130,131,132,133,134,136,137,138)

$id = 1
foreach ($r in $idRows) {
    $ws1.Cells.Item($r, 1).Value = $id
    $id = $id + 1
}

# ---------------------------------------------------------------------
# 2. Normalise the font on the "App Review" ID cells (A120:A128) so it
#    matches the lighter, non-bold font already used for the same
#    column elsewhere on the sheet.
# ---------------------------------------------------------------------
$reviewIdRange = $ws1.Range("A120:A128")
$reviewIdRange.Font.Bold = $false
$reviewIdRange.Font.Size = 10

# ---------------------------------------------------------------------
# 3. Page setup / print options for sheet 1.
# ---------------------------------------------------------------------
$ps1 = $ws1.PageSetup
$ps1.PrintGridlines = $true
$ps1.CenterHorizontally = $true
$ps1.LeftMargin = 18
$ps1.RightMargin = 18
$ps1.TopMargin = 54
$ps1.BottomMargin = 54
$ps1.HeaderMargin = 0
$ps1.FooterMargin = 0
$ps1.PaperSize = 9
$ps1.Orientation = 1
$ps1.FitToPagesWide = 1
$ps1.FitToPagesTall = 0

# ---------------------------------------------------------------------
# 4. Re-theme the "Legal" section heading font on sheet 2 (B31) from the
#    unavailable "SF Pro Display" to the workbook's Arial minor font.
# ---------------------------------------------------------------------
$legalCell = $ws2.Range("B31")
$legalCell.Font.Name = "Arial"
$legalCell.Font.ThemeFont = 1

# ---------------------------------------------------------------------
# 5. Remove the three stray blank rows at the bottom of sheet 2.
# ---------------------------------------------------------------------
$ws2.Range("A37:B39").EntireRow.Delete()
